$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add the new "amount < b" condition in column B
$ws.Range("B2").Value = "amount < b"

# Row 3: column B becomes a numeric value instead of "Travel Expenses"
$ws.Range("B3").Value = 2

# Row 4: column B becomes a numeric value instead of "Foo"
$ws.Range("B4").Value = 3.5678700000000001

# Update the active cell selection to C1
$ws.Range("C1").Select()
